# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which carry duplicate rows that must stay in sync.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 330
    "F4"  = 10480
    "F5"  = 331
    "F6"  = 952
    "F7"  = 41
    "F8"  = 1288
    "F9"  = 7260
    "F11" = 448
    "F12" = 207
    "F13" = 131
    "F17" = 686
    "F22" = 1663
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
